$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$oldText = $cellA1.Value()
$newText = $oldText.Replace("1000 Bs = 7.11 = 28492.18 pesos", "1000 Bs = 7.04 = 28257.11 pesos")
$newText = $newText.Replace("28492.18 pesos = 7.08 = 954.91 Bs", "28257.11 pesos = 7.01 = 966.77 Bs")
$cellA1.Value = $newText

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 142
$ws2.Range("O10").Value = 4012.51
$ws2.Range("N12").Value = 4029.99
$ws2.Range("O12").Value = 137.88
